$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 479.0909
$ws.Cells.Item(80, 9).Value = 266
$ws.Cells.Item(80, 10).Value = 734.8
$ws.Cells.Item(80, 11).Value = 798
$ws.Cells.Item(80, 12).Value = 2204.4
$ws.Cells.Item(80, 13).Value = 200
$ws.Cells.Item(80, 14).Value = -4200.4
$ws.Cells.Item(83, 8).Value = 479.0909
$ws.Cells.Item(83, 9).Value = 266
$ws.Cells.Item(83, 10).Value = 734.8
$ws.Cells.Item(83, 11).Value = 2394
$ws.Cells.Item(83, 12).Value = 6613.2
$ws.Cells.Item(83, 13).Value = 2598
$ws.Cells.Item(83, 14).Value = -16597.2
$ws.Cells.Item(125, 8).Value = 8268
$ws.Cells.Item(125, 9).Value = 20685.334
$ws.Cells.Item(125, 10).Value = 817.6
$ws.Cells.Item(125, 11).Value = 186168.006
$ws.Cells.Item(125, 12).Value = 7358.400000000001
$ws.Cells.Item(125, 13).Value = -183708.006
$ws.Cells.Item(125, 14).Value = -12278.4
$ws.Cells.Item(138, 8).Value = 3498.2134
$ws.Cells.Item(138, 9).Value = 2157.8235
$ws.Cells.Item(138, 10).Value = 4609.756
$ws.Cells.Item(138, 11).Value = 6473.470499999999
$ws.Cells.Item(138, 12).Value = 13829.268
$ws.Cells.Item(138, 13).Value = -1333.470499999999
$ws.Cells.Item(138, 14).Value = -24109.268

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 1020.125
$ws.Cells.Item(4, 9).Value = 1060.2
$ws.Cells.Item(4, 10).Value = 953.3333
$ws.Cells.Item(4, 11).Value = 1060.2
$ws.Cells.Item(4, 12).Value = 953.3333
$ws.Cells.Item(4, 13).Value = -944.2
$ws.Cells.Item(4, 14).Value = -1185.3333
$ws.Cells.Item(5, 8).Value = 193.75
$ws.Cells.Item(5, 9).Value = 225
$ws.Cells.Item(5, 10).Value = 162.5
$ws.Cells.Item(5, 11).Value = 225
$ws.Cells.Item(5, 12).Value = 162.5
$ws.Cells.Item(5, 13).Value = -113
$ws.Cells.Item(5, 14).Value = -386.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 193.75
$ws.Cells.Item(4, 9).Value = 225
$ws.Cells.Item(4, 10).Value = 162.5
$ws.Cells.Item(4, 11).Value = 225
$ws.Cells.Item(4, 12).Value = 162.5
$ws.Cells.Item(4, 13).Value = -110
$ws.Cells.Item(4, 14).Value = -392.5
$ws.Cells.Item(22, 8).Value = 200
$ws.Cells.Item(22, 9).Value = 200
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 200
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).ClearContents()
$ws.Cells.Item(22, 14).Value = -27
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).ClearContents()
$ws.Cells.Item(87, 14).Value = 0
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).ClearContents()
$ws.Cells.Item(90, 14).Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 529.8823
$ws.Cells.Item(22, 9).Value = 180.42857
$ws.Cells.Item(22, 10).Value = 2160.6667
$ws.Cells.Item(22, 11).Value = 180.42857
$ws.Cells.Item(22, 12).Value = 2160.6667
$ws.Cells.Item(22, 13).Value = 169.57143
$ws.Cells.Item(22, 14).Value = -2860.6667
$ws.Cells.Item(58, 8).Value = 41668020
$ws.Cells.Item(58, 9).Value = 142858100
$ws.Cells.Item(58, 10).Value = 1518.3529
$ws.Cells.Item(58, 11).Value = 142858100
$ws.Cells.Item(58, 12).Value = 1518.3529
$ws.Cells.Item(58, 13).Value = -142857897
$ws.Cells.Item(58, 14).Value = -1924.3529
$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 12).ClearContents()
$ws.Cells.Item(104, 14).Value = 0
$ws.Cells.Item(136, 8).Value = 41668020
$ws.Cells.Item(136, 9).Value = 142858100
$ws.Cells.Item(136, 10).Value = 1518.3529
$ws.Cells.Item(136, 11).Value = 428574300
$ws.Cells.Item(136, 12).Value = 4555.0587
$ws.Cells.Item(136, 13).Value = -428571750
$ws.Cells.Item(136, 14).Value = -9655.058700000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 686.73334
$ws.Cells.Item(2, 9).Value = 1319.5714
$ws.Cells.Item(2, 10).Value = 133
$ws.Cells.Item(2, 11).Value = 7917.428400000001
$ws.Cells.Item(2, 12).Value = 798
$ws.Cells.Item(2, 13).Value = -7804.428400000001
$ws.Cells.Item(2, 14).Value = -1024
$ws.Cells.Item(5, 8).Value = 792.6857
$ws.Cells.Item(5, 9).Value = 556.08
$ws.Cells.Item(5, 10).Value = 1384.2
$ws.Cells.Item(5, 11).Value = 1668.24
$ws.Cells.Item(5, 12).Value = 4152.6
$ws.Cells.Item(5, 13).Value = -1556.24
$ws.Cells.Item(5, 14).Value = -4376.6
$ws.Cells.Item(14, 8).Value = 1090.6
$ws.Cells.Item(14, 9).Value = 1090.6
$ws.Cells.Item(14, 11).Value = 3271.8
$ws.Cells.Item(14, 13).Value = -3098.8
$ws.Cells.Item(107, 8).Value = 1128.6923
$ws.Cells.Item(107, 10).Value = 1600
$ws.Cells.Item(107, 12).Value = 4800
$ws.Cells.Item(107, 14).Value = -8640
$ws.Cells.Item(122, 8).Value = 1739.2858
$ws.Cells.Item(122, 9).Value = 1200
$ws.Cells.Item(122, 10).Value = 1829.1666
$ws.Cells.Item(122, 11).Value = 10800
$ws.Cells.Item(122, 12).Value = 16462.4994
$ws.Cells.Item(122, 13).Value = -8350
$ws.Cells.Item(122, 14).Value = -21362.4994
$ws.Cells.Item(135, 8).Value = 792.6857
$ws.Cells.Item(135, 9).Value = 556.08
$ws.Cells.Item(135, 10).Value = 1384.2
$ws.Cells.Item(135, 11).Value = 5004.72
$ws.Cells.Item(135, 12).Value = 12457.8
$ws.Cells.Item(135, 13).Value = -2469.72
$ws.Cells.Item(135, 14).Value = -17527.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).ClearContents()
$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).ClearContents()
$ws.Cells.Item(81, 14).Value = 0
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).ClearContents()
$ws.Cells.Item(84, 14).Value = 0
$ws.Cells.Item(132, 8).Value = 5122.35
$ws.Cells.Item(132, 9).Value = 5707.793
$ws.Cells.Item(132, 10).Value = 3578.9092
$ws.Cells.Item(132, 11).Value = 17123.379
$ws.Cells.Item(132, 12).Value = 10736.7276
$ws.Cells.Item(132, 13).Value = -14593.379
$ws.Cells.Item(132, 14).Value = -15796.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2368.6667
$ws.Cells.Item(7, 9).Value = 1427
$ws.Cells.Item(7, 10).Value = 3122
$ws.Cells.Item(7, 11).Value = 1427
$ws.Cells.Item(7, 12).Value = 3122
$ws.Cells.Item(7, 13).Value = -1315
$ws.Cells.Item(7, 14).Value = -3346
$ws.Cells.Item(22, 8).Value = 404.9
$ws.Cells.Item(22, 9).Value = 399.66666
$ws.Cells.Item(22, 10).Value = 407.14285
$ws.Cells.Item(22, 11).Value = 399.66666
$ws.Cells.Item(22, 12).Value = 407.14285
$ws.Cells.Item(22, 13).Value = -104.66666
$ws.Cells.Item(22, 14).Value = -997.14285
$ws.Cells.Item(27, 8).Value = 404.9
$ws.Cells.Item(27, 9).Value = 399.66666
$ws.Cells.Item(27, 10).Value = 407.14285
$ws.Cells.Item(27, 11).Value = 399.66666
$ws.Cells.Item(27, 12).Value = 407.14285
$ws.Cells.Item(27, 13).Value = -292.66666
$ws.Cells.Item(27, 14).Value = -621.14285
$ws.Cells.Item(40, 8).Value = 2219.4546
$ws.Cells.Item(40, 9).Value = 1734
$ws.Cells.Item(40, 10).Value = 2802
$ws.Cells.Item(40, 11).Value = 1734
$ws.Cells.Item(40, 12).Value = 2802
$ws.Cells.Item(40, 13).Value = -1598
$ws.Cells.Item(40, 14).Value = -3074
$ws.Cells.Item(55, 8).Value = 211.73334
$ws.Cells.Item(55, 9).Value = 176
$ws.Cells.Item(55, 10).Value = 258.46155
$ws.Cells.Item(55, 11).Value = 176
$ws.Cells.Item(55, 12).Value = 258.46155
$ws.Cells.Item(55, 13).Value = -3
$ws.Cells.Item(55, 14).Value = -604.46155
$ws.Cells.Item(80, 8).Value = 48000
$ws.Cells.Item(80, 10).Value = 48000
$ws.Cells.Item(80, 12).Value = 48000
$ws.Cells.Item(80, 14).Value = -50246
$ws.Cells.Item(83, 8).Value = 48000
$ws.Cells.Item(83, 10).Value = 48000
$ws.Cells.Item(83, 12).Value = 144000
$ws.Cells.Item(83, 14).Value = -155232
$ws.Cells.Item(100, 8).Value = 2272.35
$ws.Cells.Item(100, 9).Value = 1181.8182
$ws.Cells.Item(100, 10).Value = 2686
$ws.Cells.Item(100, 11).Value = 1181.8182
$ws.Cells.Item(100, 12).Value = 2686
$ws.Cells.Item(100, 13).Value = -640.8181999999999
$ws.Cells.Item(100, 14).Value = -3768
$ws.Cells.Item(126, 8).Value = 2368.6667
$ws.Cells.Item(126, 9).Value = 1427
$ws.Cells.Item(126, 10).Value = 3122
$ws.Cells.Item(126, 11).Value = 4281
$ws.Cells.Item(126, 12).Value = 9366
$ws.Cells.Item(126, 13).Value = -1811
$ws.Cells.Item(126, 14).Value = -14306
$ws.Cells.Item(128, 8).Value = 55000
$ws.Cells.Item(128, 10).Value = 55000
$ws.Cells.Item(128, 12).Value = 55000
$ws.Cells.Item(128, 14).Value = -64960
$ws.Cells.Item(136, 8).Value = 5385.875
$ws.Cells.Item(136, 9).Value = 2019
$ws.Cells.Item(136, 10).Value = 12793
$ws.Cells.Item(136, 11).Value = 6057
$ws.Cells.Item(136, 12).Value = 38379
$ws.Cells.Item(136, 13).Value = -3507
$ws.Cells.Item(136, 14).Value = -43479

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 319.46155
$ws.Cells.Item(107, 9).Value = 241.33333
$ws.Cells.Item(107, 10).Value = 495.25
$ws.Cells.Item(107, 11).Value = 723.99999
$ws.Cells.Item(107, 12).Value = 1485.75
$ws.Cells.Item(107, 13).Value = 1196.00001
$ws.Cells.Item(107, 14).Value = -5325.75
$ws.Cells.Item(126, 8).Value = 817.6667
$ws.Cells.Item(126, 9).Value = 402
$ws.Cells.Item(126, 10).Value = 1233.3334
$ws.Cells.Item(126, 11).Value = 1206
$ws.Cells.Item(126, 12).Value = 3700.0002
$ws.Cells.Item(126, 13).Value = 1264
$ws.Cells.Item(126, 14).Value = -8640.0002
$ws.Cells.Item(132, 8).Value = 3154.75
$ws.Cells.Item(132, 9).Value = 3708.26
$ws.Cells.Item(132, 10).Value = 2090.3076
$ws.Cells.Item(132, 11).Value = 11124.78
$ws.Cells.Item(132, 12).Value = 6270.9228
$ws.Cells.Item(132, 13).Value = -8594.780000000001
$ws.Cells.Item(132, 14).Value = -11330.9228
